# Split / re-flow the "inflow method" description strings: the
# Extreme Low Flow Method (col A) and Strategy to Stabilize Lake Levels
# (col F) text cells on the ExtremeFlows sheet pick up an extra internal
# space where the string had been split and rejoined. Also move the
# active selection on the ExtremeFlows sheet from F8 to F7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtremeFlows")

# Column A - Extreme Low Flow Method labels
$ws.Range("A3").Value = "B. Collaborator choices in immersive  modeling sessions"
$ws.Range("A4").Value = "C. Low Lake Powell releases + gains  through Grand Canyon"
$ws.Range("A5").Value = "D. 85%, 65%, and 50% of 2000 to 2018  average flow"
$ws.Range("A6").Value = "E. Reclamation's Post 2026  web tool"
$ws.Range("A7").Value = "F. Lowest consecutive flows in Reclamation's  ensembles and traces"

# Column F - Strategy to Stabilize Lake Levels descriptions
$ws.Range("F2").Value = "Cap depletions to 10-year lookback  period of flow."
$ws.Range("F3").Value = "Divide inflow; Users consume and conserve  within their account balance."
$ws.Range("F4").Value = "Rule curve; Consumption equals or less  than inflow minus evaporation."
$ws.Range("F5").Value = "Release 95% of regulated  inflow."
$ws.Range("F6").Value = "Release to prevent drawdown  to 3,490 feet."

# Move the active cell selection from F8 to F7
$ws.Range("F7").Select()
